$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-06-14 Saturday" "2025-06-15 Sunday"

Replace-Text "834÷5=" "621÷3="
Replace-Text "613÷4=" "630÷4="
Replace-Text "930÷3=" "951÷3="
Replace-Text "555÷4=" "678÷4="
Replace-Text "349÷3=" "955÷7="

Replace-Text "965÷3=" "404÷3="
Replace-Text "395÷5=" "914÷9="
Replace-Text "135÷4=" "280÷3="
Replace-Text "753÷3=" "657÷6="
Replace-Text "216÷9=" "921÷3="

Replace-Text "704÷7=" "941÷8="
Replace-Text "959÷2=" "103÷4="
Replace-Text "157÷8=" "377÷2="
Replace-Text "120÷4=" "180÷4="
Replace-Text "276÷9=" "235÷2="

Replace-Text "646÷3=" "302÷8="
Replace-Text "858÷3=" "679÷8="
Replace-Text "754÷3=" "334÷4="
Replace-Text "269÷6=" "990÷2="
Replace-Text "699÷7=" "936÷8="

Replace-Text "969÷3=" "648÷7="
Replace-Text "284÷3=" "794÷2="
Replace-Text "366÷3=" "226÷3="
Replace-Text "626÷4=" "636÷7="
Replace-Text "707÷7=" "895÷5="
